$d = $word.ActiveDocument
$r = $d.Content

# Append a new paragraph after the existing "Halo World" paragraph
# containing "This is some other data".
$r.InsertParagraphAfter()
$r.InsertAfter("This is some other data")

# Append another new paragraph after that one containing
# "And some more " (note trailing space).
$r.InsertParagraphAfter()
$r.InsertAfter("And some more ")
